$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header renames -------------------------------------------------

$ws.Range("A1").Value = "button_alertActions_class"

$ws.Range("D1").Value = "div_testRuns_internalRoleCellName"
$ws.Range("E1").Value = "div_testRuns_internalRoleCellName_1"
$ws.Range("F1").Value = "div_testRuns_internalRoleCellName_10"
$ws.Range("G1").Value = "div_testRuns_internalRoleCellName_2"
$ws.Range("H1").Value = "div_testRuns_internalRoleCellName_3"
$ws.Range("I1").Value = "div_testRuns_internalRoleCellName_4"
$ws.Range("J1").Value = "div_testRuns_internalRoleCellName_5"
$ws.Range("K1").Value = "div_testRuns_internalRoleCellName_6"
$ws.Range("L1").Value = "div_testRuns_internalRoleCellName_7"
$ws.Range("M1").Value = "div_testRuns_internalRoleCellName_8"
$ws.Range("N1").Value = "div_testRuns_internalRoleCellName_9"

$ws.Range("P1").Value = "header_testRunTitles_internalRoleHeadingName"
$ws.Range("Q1").Value = "header_testRunTitles_internalRoleHeadingName_1"
$ws.Range("R1").Value = "header_testRunTitles_internalRoleHeadingName_2"

$ws.Range("AE1").Value = "link_executionLinks_executions_id"
$ws.Range("AF1").Value = "link_executionLinks_executions_id_1"
$ws.Range("AG1").Value = "link_executionLinks_internalRoleLinkName"
$ws.Range("AH1").Value = "link_executionLinks_internalRoleLinkName_1"
$ws.Range("AI1").Value = "link_executionLinks_project_id"
$ws.Range("AJ1").Value = "link_executionLinks_project_id_1"
$ws.Range("AK1").Value = "link_executionLinks_team_id"
$ws.Range("AL1").Value = "link_executionLinks_team_id_1"

$ws.Range("AM1").Value = "link_projectLinks_internalRoleLinkName"
$ws.Range("AN1").Value = "link_projectLinks_project_id"
$ws.Range("AO1").Value = "link_projectLinks_team_id"
$ws.Range("AP1").Value = "link_projectLinks_test_project_id"
$ws.Range("AQ1").Value = "link_projectLinks_trNthChild"

$ws.Range("AR1").Value = "link_testRunLinks_plan_id"
$ws.Range("AS1").Value = "link_testRunLinks_plan_id_1"
$ws.Range("AT1").Value = "link_testRunLinks_project_id"
$ws.Range("AU1").Value = "link_testRunLinks_project_id_1"
$ws.Range("AV1").Value = "link_testRunLinks_team_id"
$ws.Range("AW1").Value = "link_testRunLinks_team_id_1"

# --- Row 2 data path updates ----------------------------------------------

$ws.Range("B2").Value = "Data Files/AI-Generated/Common/scheduleTestRunAndConfigureEnvironment-test-data"
$ws.Range("C2").Value = "Data Files/AI-Generated/Common/scheduleTestRunAndConfigureEnvironment-test-data"

# --- Re-autofit the columns whose header/content text changed -------------
# (mirrors Excel auto-resizing columns after the cell content was edited)

$ws.Range("B1:N1").EntireColumn.AutoFit()
$ws.Range("P1:R1").EntireColumn.AutoFit()
$ws.Range("AE1:AW1").EntireColumn.AutoFit()
